# Add a new "Clustering Channels" column to the "beads" worksheet, letting
# the user specify which channels to use when clustering the beads data.

$wb = $excel.ActiveWorkbook
$beads = $wb.Worksheets.Item("beads")

# Match the formatting of the other header / data cells in the row:
# E1 should look like D1 (bold header style), E2 like D2 (default style).
$beads.Range("D1").Copy()
$beads.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$beads.Range("E1").Value = "Clustering Channels"

$beads.Range("D2").Copy()
$beads.Range("E2").PasteSpecial(-4122)  # xlPasteFormats
$beads.Range("E2").Value = "FL1-H, FL2-H, FL3-H"

$beads.Application.CutCopyMode = $false

# Give the new column a sensible width, matching the other bestFit columns.
$beads.Columns.Item(5).ColumnWidth = 18

# The "beads" sheet (and cell D6 on it) is now the active selection instead
# of cell E3 on the "cells" sheet.
$beads.Activate() | Out-Null
$beads.Range("D6").Select() | Out-Null
